$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column = 想去人数 / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14886
$ws1.Range("F3").Value = 18551
$ws1.Range("F22").Value = 7695
$ws1.Range("F29").Value = 102
$ws1.Range("F34").Value = 5312
$ws1.Range("F36").Value = 39

# Sheet "全部类型" (same rows updated, mirrored data)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14886
$ws4.Range("F3").Value = 18551
$ws4.Range("F23").Value = 7695
$ws4.Range("F32").Value = 102
$ws4.Range("F37").Value = 5312
$ws4.Range("F39").Value = 39
